$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1,2,3,5,8,9,10 -> -0.252188276307633
$rowsA = @(1,2,3,5,8,9,10)
foreach ($r in $rowsA) {
    $ws.Range("A$r`:J$r").Value = -0.252188276307633
}

# Row 6 -> -0.2510504393851208
$ws.Range("A6:J6").Value = -0.2510504393851208

# Row 7 -> -1.111953237469457
$ws.Range("A7:J7").Value = -1.111953237469457
